$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before row 812, shifting existing rows 812+ down to 814+
$ws.Rows("812:813").Insert()

# Row 812: based on the original row 812 (now at row 814) template, with updated values
$ws.Range("A812").Value = 3
$ws.Range("B812").Value = "Femacal de La Calera"
$ws.Range("C812").Value = "Coquimbo"
$ws.Range("D812").Value = 44931
$ws.Range("D812").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E812").Value = 5
$ws.Range("F812").Value = "Fruta"
$ws.Range("G812").Value = 100108
$ws.Range("H812").Value = "Tropicales y subtropicales"
$ws.Range("I812").Value = 100108005
$ws.Range("J812").Value = "Piña"
$ws.Range("K812").Value = "Caramelo"
$ws.Range("L812").Value = "Primera"
$ws.Range("M812").Value = 162
$ws.Range("N812").Value = 21000
$ws.Range("O812").Value = 21000
$ws.Range("P812").Value = 21000
$ws.Range("Q812").Value = "$/caja 12 unidades"
$ws.Range("R812").Value = "Ecuador"
$ws.Range("S812").Value = 1750
$ws.Range("T812").Value = 12

# Row 813: based on original row 813 template, with updated values
$ws.Range("A813").Value = 3
$ws.Range("B813").Value = "Femacal de La Calera"
$ws.Range("C813").Value = "Coquimbo"
$ws.Range("D813").Value = 44931
$ws.Range("D813").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E813").Value = 5
$ws.Range("F813").Value = "Fruta"
$ws.Range("G813").Value = 100108
$ws.Range("H813").Value = "Tropicales y subtropicales"
$ws.Range("I813").Value = 100108005
$ws.Range("J813").Value = "Piña"
$ws.Range("K813").Value = "Caramelo"
$ws.Range("L813").Value = "Segunda"
$ws.Range("M813").Value = 162
$ws.Range("N813").Value = 21000
$ws.Range("O813").Value = 21000
$ws.Range("P813").Value = 21000
$ws.Range("Q813").Value = "$/caja 14 unidades"
$ws.Range("R813").Value = "Ecuador"
$ws.Range("S813").Value = 1500
$ws.Range("T813").Value = 14
